# Updated cryptos list on Sat Sep  7 16:42:02 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, [string]$value) {
    # Force the cell to be treated as text so numeric-looking strings
    # (e.g. "0.529", "21.79") are not silently converted to numbers,
    # then restore the default ("Normal") style so no stray formatting
    # is left behind on the cell.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "54.171.47"
$ws.Range("E2").Value = "  +0.23%  "

Set-TextCell $ws.Range("D3") "2.262.50"
$ws.Range("E3").Value = "  -0.28%  "

$ws.Range("E4").Value = "  +0.21%  "

Set-TextCell $ws.Range("D5") "496.19"
$ws.Range("E5").Value = "  +1.58%  "

Set-TextCell $ws.Range("D6") "127.78"
$ws.Range("E6").Value = "  +0.26%  "

$ws.Range("E7").Value = "  +0.02%  "

Set-TextCell $ws.Range("D8") "0.529"
$ws.Range("E8").Value = "  +1.04%  "

Set-TextCell $ws.Range("D9") "2.287.49"
$ws.Range("E9").Value = "  +0.60%  "

Set-TextCell $ws.Range("D10") "0.0948"
$ws.Range("E10").Value = "  +2.37%  "

$ws.Range("E11").Value = "  +2.19%  "

$ws.Range("E12").Value = "  +2.68%  "

$ws.Range("E13").Value = "  -3.17%  "

Set-TextCell $ws.Range("D14") "2.701.74"
$ws.Range("E14").Value = "  +0.95%  "

Set-TextCell $ws.Range("D15") "21.79"
$ws.Range("E15").Value = "  +3.24%  "

Set-TextCell $ws.Range("D16") "54.499.63"
$ws.Range("E16").Value = "  +0.88%  "

Set-TextCell $ws.Range("D17") "0.0000129"
$ws.Range("E17").Value = "  +0.60%  "

Set-TextCell $ws.Range("D18") "2.327.23"
$ws.Range("E18").Value = "  +1.85%  "

Set-TextCell $ws.Range("D19") "10.07"
$ws.Range("E19").Value = "  +4.60%  "

Set-TextCell $ws.Range("D20") "4.11"
$ws.Range("E20").Value = "  +3.19%  "

Set-TextCell $ws.Range("D21") "6.48"
$ws.Range("E21").Value = "  +5.52%  "

Set-TextCell $ws.Range("D22") "302.41"
$ws.Range("E22").Value = "  -0.12%  "

$ws.Range("E23").Value = "  -0.04%  "

$ws.Range("E24").Value = "  -1.96%  "

Set-TextCell $ws.Range("D25") "62.70"
$ws.Range("E25").Value = "  -1.99%  "

$ws.Range("E26").Value = "  +0.13%  "

Set-TextCell $ws.Range("D27") "0.374"
$ws.Range("E27").Value = "  +2.12%  "

Set-TextCell $ws.Range("D28") "2.405.32"
$ws.Range("E28").Value = "  +0.79%  "

Set-TextCell $ws.Range("D29") "0.151"
$ws.Range("E29").Value = "  +5.09%  "

Set-TextCell $ws.Range("D30") "7.09"
$ws.Range("E30").Value = "  +0.44%  "

Set-TextCell $ws.Range("D31") "169.34"
$ws.Range("E31").Value = "  -0.11%  "

Set-TextCell $ws.Range("D32") "0.0₃0691"
$ws.Range("E32").Value = "  -0.62%  "

$ws.Range("E33").Value = "  +0.23%  "

Set-TextCell $ws.Range("D34") "5.89"
$ws.Range("E34").Value = "  +2.23%  "

$ws.Range("E35").Value = "  +0.08%  "

Set-TextCell $ws.Range("D36") "0.992"
$ws.Range("E36").Value = "  -0.69%  "

Set-TextCell $ws.Range("D37") "1.07"
$ws.Range("E37").Value = "  +0.60%  "

Set-TextCell $ws.Range("D38") "17.64"
$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("E39").Value = "  +2.67%  "

Set-TextCell $ws.Range("D40") "0.868"
$ws.Range("E40").Value = "  +3.92%  "

$ws.Range("E41").Value = "  +2.30%  "

Set-TextCell $ws.Range("D42") "35.45"
$ws.Range("E42").Value = "  -1.08%  "

$ws.Range("E43").Value = "  +2.71%  "

Set-TextCell $ws.Range("D44") "0.375"
$ws.Range("E44").Value = "  +2.19%  "

Set-TextCell $ws.Range("D45") "3.35"
$ws.Range("E45").Value = "  +1.15%  "

Set-TextCell $ws.Range("D46") "128.16"
$ws.Range("E46").Value = "  +4.65%  "

Set-TextCell $ws.Range("D47") "4.80"
$ws.Range("E47").Value = "  +2.20%  "

Set-TextCell $ws.Range("D48") "0.0890"
$ws.Range("E48").Value = "  +1.25%  "

Set-TextCell $ws.Range("D49") "0.544"
$ws.Range("E49").Value = "  +0.62%  "

Set-TextCell $ws.Range("D50") "240.53"
$ws.Range("E50").Value = "  +0.52%  "

Set-TextCell $ws.Range("D51") "0.0485"
